# Update cell A1 on the active sheet ("Feuil1") from "Donnée A1 " to "modif 10h30".
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "modif 10h30"

# Restore the selection/active cell to A1 (the file's selection moves off the
# stray "F5" leftover back onto the edited cell).
[void]$ws.Range("A1").Select()
